$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.328.49'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '1.947.91'
$ws.Range('E3').Value = '  -3.81%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.611'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '57.50'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -8.52%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.367'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '55.46'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0819'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.31%  '
$ws.Range('E12').Value = '  -0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.826'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.09%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.238.55'
$ws.Range('E14').Value = '  -3.34%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.25'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -9.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.32%  '
$ws.Range('D18').Value = '1.961.28'
$ws.Range('E18').Value = '  -3.11%  '
$ws.Range('D19').Value = '36.261.25'
$ws.Range('E19').Value = '  -1.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.23%  '
$ws.Range('D21').Value = '0.0₃0867'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.57%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.40'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.33'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.117'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -14.68%  '
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.14'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.69'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.42%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0630'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.79'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.14'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -10.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -12.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0980'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.73%  '
$ws.Range('E42').Value = '  -2.89%  '
$ws.Range('E43').Value = '  -6.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0209'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -8.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.03'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.35'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.21%  '
$ws.Range('D48').Value = '1.337.89'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.81%  '
